# "solved binary tree level order traversal"
# Append a new tracker row (row 21) for the "Binary Tree Level Order
# Traversal" LeetCode problem, mirroring the existing table layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fill in the new row's values first, in left-to-right column order,
#    so new shared-string entries are created in the same order the
#    original author would have typed them (number, difficulty,
#    question, url, approach, efficient-O, highlight).
$ws.Range("A21").Value = 102
$ws.Range("B21").Value = "Medium"
$ws.Range("C21").Value = "Binary Tree Level Order Traversal"
$ws.Range("D21").Value = "https://shorturl.at/89NnQ"
$ws.Range("E21").Value = "Queue"
$ws.Range("F21").Value = "O(n)"
$ws.Range("G21").Value = "Use DFS using queues and use a counter variable (size) to get the level arrays inside the res array."

# 2. Turn the url text in D21 into a real hyperlink.
$ws.Hyperlinks.Add($ws.Range("D21"), "https://shorturl.at/89NnQ") | Out-Null

# 3. Match the formatting used by the rest of the "Easy"-style rows
#    (e.g. row 20) by copying its formats down onto the new row.
$ws.Range("A20:G20").Copy() | Out-Null
$ws.Range("A21:G21").PasteSpecial(-4122) | Out-Null

# 4. Leave the selection where the author's cursor ended up after
#    entering the new row.
$ws.Range("C23").Select() | Out-Null
